# Apply the commit's changes:
#  1. Slide 6's table is re-styled from the custom "Integral" table style
#     to the built-in PowerPoint table style {23760D92-FB49-4159-A0F2-BD9C6E919199}.
#  2. The presentation's theme colour scheme is changed from the custom
#     "Integral" palette to the default "Office" palette (dk1, lt1, dk2, lt2,
#     accent1-6, hlink, folHlink - the 12 DrawingML theme colours).

$p = $ppt.ActivePresentation

# --- 1. Re-apply the table style on the table located on slide 6 ---------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{23760D92-FB49-4159-A0F2-BD9C6E919199}", $true)
    }
}

# --- 2. Swap the theme colour scheme to the default Office palette -------
# ThemeColorScheme index -> (scheme slot, new RGB as a COM BGR-encoded long)
#   1  dk1       000000 -> 0
#   2  lt1       FFFFFF -> 16777215
#   3  dk2       44546A -> 6968388
#   4  lt2       E7E6E6 -> 15132391
#   5  accent1   5B9BD5 -> 13998939
#   6  accent2   ED7D31 -> 3243501
#   7  accent3   A5A5A5 -> 10855845
#   8  accent4   FFC000 -> 49407
#   9  accent5   4472C4 -> 12874308
#   10 accent6   70AD47 -> 4697456
#   11 hlink     0563C1 -> 12673797
#   12 folHlink  954F72 -> 7491477
$officeColors = @(0, 16777215, 6968388, 15132391, 13998939, 3243501, 10855845, 49407, 12874308, 4697456, 12673797, 7491477)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
